$wb = $excel.ActiveWorkbook

# Find the last sheet (most recent week) to use as a style/layout template
$lastIndex = $wb.Worksheets.Count
$srcSheet = $wb.Worksheets.Item($lastIndex)

# Duplicate it (copy lands immediately after the source) so the new sheet
# inherits the exact same formatting (bold/bordered/centered header style)
$srcSheet.Copy($null, $srcSheet)

$newSheet = $wb.Worksheets.Item($lastIndex + 1)
$newSheet.Name = "magapoke_2025-12-17"

# Header row (rank / title) is already correct on the copied sheet;
# column A (rank 1..100) is also already correct since it is the same
# sequential list every week. Only the title column (B) changes.
$titles = @(
    'ブルーロック',
    '東京卍リベンジャーズ',
    '信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！',
    'ベイビーステップ',
    'みいちゃんと山田さん',
    'ギルティサークル',
    'ドラハチ',
    'ガチアクタ',
    '島耕作',
    'イレギュラーズ',
    '薫る花は凛と咲く',
    '黄昏町プリズナーズ',
    '十字架のろくにん',
    '黒猫と魔女の教室',
    'ハードワーカー中田',
    'WIND BREAKER',
    '魔女と傭兵',
    'となりの黒川さん',
    '魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～',
    '転生したら第七王子だったので、気ままに魔術を極めます',
    '君が僕らを悪魔と呼んだ頃',
    '異世界ウォーキング',
    '蒼く染めろ',
    '味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す',
    'ハナバス　苔石花江のバスケ論',
    'K-9~警視庁公安部公安第9課異能対策係~',
    'アルキメデスの大戦',
    'せいぶつ部の田辺くん',
    '南海トラフ巨大地震',
    '辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜',
    'グラぱらっ！',
    '転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～',
    '幼馴染とはラブコメにならない',
    'FAIRY TAIL 100 YEARS QUEST',
    'ひゃくえむ。',
    'アオバノバスケ',
    'ナキナギ',
    '愛妻の裏アカ',
    '時々ボソッとロシア語でデレる隣のアーリャさん',
    'いじめるヤバイ奴',
    'さわらないで小手指くん',
    '食糧人類-Starving Anonymous-',
    '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～',
    '屋根の下のアルテミス',
    '不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～',
    'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！',
    '不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～',
    'おやすみ ふみさん',
    '念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！',
    'デッドアカウント',
    '普通の本はありません！',
    '春くらり',
    '降り積もれ孤独な死よ',
    '東京卍リベンジャーズ～場地圭介からの手紙～',
    '追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜',
    'ジュミドロ',
    'デスティニーラバーズ',
    '我間乱 ―修羅―',
    'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜',
    '最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～',
    '阿武ノーマル',
    '君が監督！',
    '冰剣の魔術師が世界を統べる〜世界最強の魔術師である少年は、魔術学院に入学する〜',
    'お嬢様の僕',
    '劣等人の魔剣使い　スキルボードを駆使して最強に至る',
    'MYS',
    'Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜',
    'Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。',
    'お願い、脱がシて。',
    '卒業アルバムの彼女たち',
    '恋ニ非ズ',
    'リスナーに騙されてダンジョンの最下層から脱出RTAすることになった',
    'ストーカー行為がバレて人生終了男',
    'インフェクション',
    '英雄と魔女の転生ラブコメ',
    '田んぼで拾った女騎士、田舎で俺の嫁だと思われている',
    '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～',
    '金田一少年の事件簿外伝 犯人たちの事件簿',
    'GALAXIAS',
    '可愛いだけじゃない式守さん',
    'ヒロインは絶望しました。',
    '勇者と呼ばれた後に　―そして無双男は家族を創る―',
    'ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～',
    'ぼくのアデリア',
    '追放されなかった男　～二度目の人生は土下座から始まりました～',
    'ともだちづくり',
    '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～',
    'なれの果ての僕ら',
    'それがメイドのカンナです',
    '東京ネオンスキャンダル',
    '陽子さん、すがりよる。',
    'この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～',
    'イジらないで、長瀞さん',
    '魁の花巫女',
    '剣帝学院の魔眼賢者',
    '人間消失',
    '復讐の教科書',
    '母という呪縛 娘という牢獄',
    '四刀流の最強配信者～やり込んだVRゲームの設定が現実世界に反映されたので、廃止予定だった戦闘職で無双します～',
    'ウイニング パス'
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 2).Value = $titles[$i]
}

# Make the new weekly sheet the active one, mirroring a freshly-added tab
$newSheet.Activate()
